$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A2').Value = 38
$ws.Range('B2').Value = '$\eta_{q}$'
$ws.Range('C2').Value = 0.0003120172204467836
$ws.Range('A3').Value = 67
$ws.Range('B3').Value = '$(F_{r}^{\text{SCF}})_{3}$'
$ws.Range('C3').Value = 0.0001778814710337454
$ws.Range('A4').Value = 35
$ws.Range('B4').Value = '$F_{q}^{\text{SCF}}$'
$ws.Range('C4').Value = 0.0001668519649732861
$ws.Range('A5').Value = 26
$ws.Range('B5').Value = 'typ_3'
$ws.Range('C5').Value = 0.0001472818762419537
$ws.Range('A6').Value = 37
$ws.Range('B6').Value = '$F_{q}$'
$ws.Range('C6').Value = 0.0000774853219261417
$ws.Range('A7').Value = 23
$ws.Range('B7').Value = 'typ_0'
$ws.Range('C7').Value = 0.00007060314944748904
$ws.Range('A8').Value = 3
$ws.Range('B8').Value = 'h$_{p}^{3}$'
$ws.Range('C8').Value = 0.00005401107090161937
$ws.Range('A9').Value = 43
$ws.Range('B9').Value = '$F_{s}^{\text{SCF}}$'
$ws.Range('C9').Value = 0.00005009582467455623
$ws.Range('A10').Value = 46
$ws.Range('B10').Value = '$\eta_{s}$'
$ws.Range('C10').Value = 0.00004558262841241548
$ws.Range('A11').Value = 24
$ws.Range('B11').Value = 'typ_1'
$ws.Range('C11').Value = 0.00004318485549334553
$ws.Range('A12').Value = 76
$ws.Range('B12').Value = '$\langle ss \vert ss \rangle$'
$ws.Range('C12').Value = 0.00003236620385005609
$ws.Range('A13').Value = 12
$ws.Range('B13').Value = 'h$_{q}$'
$ws.Range('C13').Value = 0.00003194477308960642
$ws.Range('A14').Value = 45
$ws.Range('B14').Value = '$F_{s}$'
$ws.Range('C14').Value = 0.00003128022965747469
$ws.Range('A15').Value = 25
$ws.Range('B15').Value = 'typ_2'
$ws.Range('C15').Value = 0.00002812836244031945
$ws.Range('A16').Value = 74
$ws.Range('B16').Value = '$\langle qq \vert qq \rangle$'
$ws.Range('C16').Value = 0.00002437297645613405
$ws.Range('A17').Value = 94
$ws.Range('B17').Value = '$(\langle pq \vert qp \rangle)_{2}$'
$ws.Range('C17').Value = 0.00002428048765225232
$ws.Range('A18').Value = 92
$ws.Range('B18').Value = '$(\langle rr \vert rr \rangle)_{2}$'
$ws.Range('C18').Value = 0.00002154245651384006
$ws.Range('A19').Value = 91
$ws.Range('B19').Value = '$(\langle pp \vert pp \rangle)_{2}$'
$ws.Range('C19').Value = 0.00002057768493538476
$ws.Range('A20').Value = 99
$ws.Range('B20').Value = '$(\langle pp \vert pp \rangle)_{3}$'
$ws.Range('C20').Value = 0.00001965532885650237
$ws.Range('A21').Value = 93
$ws.Range('B21').Value = '$(\langle pq \vert pq \rangle)_{2}$'
$ws.Range('C21').Value = 0.00001817392497424592
$ws.Range('A22').Value = 22
$ws.Range('B22').Value = 'h$_{s}$'
$ws.Range('C22').Value = 0.00001799234872952975
$ws.Range('A23').Value = 0
$ws.Range('B23').Value = 'h$_{p}^{0}$'
$ws.Range('C23').Value = 0.00001426192432767322
$ws.Range('A24').Value = 1
$ws.Range('B24').Value = 'h$_{p}^{1}$'
$ws.Range('C24').Value = 0.0000138237298676159
$ws.Range('A25').Value = 96
$ws.Range('B25').Value = '$(\langle rs \vert sr \rangle)_{2}$'
$ws.Range('C25').Value = 0.00001330714091716932
$ws.Range('A26').Value = 95
$ws.Range('B26').Value = '$(\langle rs\vert rs \rangle)_{2}$'
$ws.Range('C26').Value = 0.00001286163063139556
$ws.Range('A27').Value = 78
$ws.Range('B27').Value = '$(\langle pq \vert qp \rangle)_{0}$'
$ws.Range('C27').Value = 0.00001230002396985267
$ws.Range('A28').Value = 51
$ws.Range('B28').Value = '$(F_{r}^{\text{SCF}})_{1}$'
$ws.Range('C28').Value = 0.00001083907340948885
$ws.Range('A29').Value = 57
$ws.Range('B29').Value = '$(F_{p})_{2}$'
$ws.Range('C29').Value = 0.0000103036158299805
$ws.Range('A30').Value = 17
$ws.Range('B30').Value = 'h$_{r}^{3}$'
$ws.Range('C30').Value = 0.000009369754075793779
$ws.Range('A31').Value = 86
$ws.Range('B31').Value = '$(\langle pq \vert qp \rangle)_{1}$'
$ws.Range('C31').Value = 0.000009278274433767798
$ws.Range('A32').Value = 100
$ws.Range('B32').Value = '$(\langle rr \vert rr \rangle)_{3}$'
$ws.Range('C32').Value = 0.000008884588843363075
$ws.Range('A33').Value = 11
$ws.Range('B33').Value = 'h$_{pr}^{3}$'
$ws.Range('C33').Value = 0.000008141327970100411
$ws.Range('A34').Value = 16
$ws.Range('B34').Value = 'h$_{r}^{2}$'
$ws.Range('C34').Value = 0.00000809640560594986
$ws.Range('A35').Value = 15
$ws.Range('B35').Value = 'h$_{r}^{1}$'
$ws.Range('C35').Value = 0.00000794883774775833
$ws.Range('A36').Value = 61
$ws.Range('B36').Value = '$(F_{r})_{2}$'
$ws.Range('C36').Value = 0.0000078281936721298
$ws.Range('A37').Value = 70
$ws.Range('B37').Value = '$(\eta_{r})_{3}$'
$ws.Range('C37').Value = 0.00000775031351713677
$ws.Range('A38').Value = 59
$ws.Range('B38').Value = '$(F_{r}^{\text{SCF}})_{2}$'
$ws.Range('C38').Value = 0.000007574674842897949
$ws.Range('A39').Value = 54
$ws.Range('B39').Value = '$(\eta_{r})_{1}$'
$ws.Range('C39').Value = 0.00000733865611703463
$ws.Range('A40').Value = 55
$ws.Range('B40').Value = '$(F_{p}^{\text{SCF}})_{2}$'
$ws.Range('C40').Value = 0.000007138719100143492
$ws.Range('A41').Value = 101
$ws.Range('B41').Value = '$(\langle pq \vert pq \rangle)_{3}$'
$ws.Range('C41').Value = 0.000007098883784052624
$ws.Range('A42').Value = 102
$ws.Range('B42').Value = '$(\langle pq \vert qp \rangle)_{3}$'
$ws.Range('C42').Value = 0.000006707839937621344
$ws.Range('A43').Value = 39
$ws.Range('B43').Value = '$(F_{r}^{\text{SCF}})_{0}$'
$ws.Range('C43').Value = 0.000006247818660499038
$ws.Range('A44').Value = 63
$ws.Range('B44').Value = '$(F_{p}^{\text{SCF}})_{3}$'
$ws.Range('C44').Value = 0.000004435591140395095
$ws.Range('A45').Value = 10
$ws.Range('B45').Value = 'h$_{pr}^{2}$'
$ws.Range('C45').Value = 0.000004414613521288797
$ws.Range('A46').Value = 2
$ws.Range('B46').Value = 'h$_{p}^{2}$'
$ws.Range('C46').Value = 0.000004359525421249219
$ws.Range('A47').Value = 69
$ws.Range('B47').Value = '$(F_{r})_{3}$'
$ws.Range('C47').Value = 0.000004240598456688464
$ws.Range('A48').Value = 97
$ws.Range('B48').Value = '$(\langle pq \vert rs \rangle)_{3}$'
$ws.Range('C48').Value = 0.00000413665581098624
$ws.Range('A49').Value = 104
$ws.Range('B49').Value = '$(\langle rs \vert sr \rangle)_{3}$'
$ws.Range('C49').Value = 0.000003937358021773291
$ws.Range('A50').Value = 47
$ws.Range('B50').Value = '$(F_{p}^{\text{SCF}})_{1}$'
$ws.Range('C50').Value = 0.000003706648974165122
$ws.Range('A51').Value = 62
$ws.Range('B51').Value = '$(\eta_{r})_{2}$'
$ws.Range('C51').Value = 0.000003563339621298202
$ws.Range('A52').Value = 31
$ws.Range('B52').Value = '$(F_{p}^{\text{SCF}})_{0}$'
$ws.Range('C52').Value = 0.000003545605974376752
$ws.Range('A53').Value = 13
$ws.Range('B53').Value = 'h$_{qs}$'
$ws.Range('C53').Value = 0.000003467056567448506
$ws.Range('A54').Value = 75
$ws.Range('B54').Value = '$(\langle rr \vert rr \rangle)_{0}$'
$ws.Range('C54').Value = 0.000003431753266815229
$ws.Range('A55').Value = 30
$ws.Range('B55').Value = 'FA$_{qs}$'
$ws.Range('C55').Value = 0.000003267765764094631
$ws.Range('A56').Value = 77
$ws.Range('B56').Value = '$(\langle pq \vert pq \rangle)_{0}$'
$ws.Range('C56').Value = 0.000002869427316362245
$ws.Range('A57').Value = 42
$ws.Range('B57').Value = '$(\eta_{r})_{0}$'
$ws.Range('C57').Value = 0.000002711234104729075
$ws.Range('A58').Value = 88
$ws.Range('B58').Value = '$(\langle rs \vert sr \rangle)_{1}$'
$ws.Range('C58').Value = 0.00000262975976143675
$ws.Range('A59').Value = 5
$ws.Range('B59').Value = 'h$_{pq}^{1}$'
$ws.Range('C59').Value = 0.000002603771348045123
$ws.Range('A60').Value = 8
$ws.Range('B60').Value = 'h$_{pr}^{0}$'
$ws.Range('C60').Value = 0.00000257900315329351
$ws.Range('A61').Value = 49
$ws.Range('B61').Value = '$(F_{p})_{1}$'
$ws.Range('C61').Value = 0.000002460173160205581
$ws.Range('A62').Value = 29
$ws.Range('B62').Value = 'FI$_{qs}$'
$ws.Range('C62').Value = 0.000002323367225995174
$ws.Range('A63').Value = 14
$ws.Range('B63').Value = 'h$_{r}^{0}$'
$ws.Range('C63').Value = 0.00000231888980937027
$ws.Range('A64').Value = 9
$ws.Range('B64').Value = 'h$_{pr}^{1}$'
$ws.Range('C64').Value = 0.000002298812190535319
$ws.Range('A65').Value = 53
$ws.Range('B65').Value = '$(F_{r})_{1}$'
$ws.Range('C65').Value = 0.000002294214471601889
$ws.Range('A66').Value = 79
$ws.Range('B66').Value = '$(\langle rs\vert rs \rangle)_{0}$'
$ws.Range('C66').Value = 0.000002200308836325464
$ws.Range('A67').Value = 73
$ws.Range('B67').Value = '$(\langle pp \vert pp \rangle)_{0}$'
$ws.Range('C67').Value = 0.000001927698869691738
$ws.Range('A68').Value = 80
$ws.Range('B68').Value = '$(\langle rs \vert sr \rangle)_{0}$'
$ws.Range('C68').Value = 0.000001694889497391638
$ws.Range('A69').Value = 71
$ws.Range('B69').Value = '$(\langle pq \vert rs \rangle)_{0}$'
$ws.Range('C69').Value = 0.000001389867104354688
$ws.Range('A70').Value = 84
$ws.Range('B70').Value = '$(\langle rr \vert rr \rangle)_{1}$'
$ws.Range('C70').Value = 0.000001312886564480988
$ws.Range('A71').Value = 19
$ws.Range('B71').Value = 'h$_{rs}^{1}$'
$ws.Range('C71').Value = 0.000001206476478763909
$ws.Range('A72').Value = 18
$ws.Range('B72').Value = 'h$_{rs}^{0}$'
$ws.Range('C72').Value = 0.000001163262892350264
$ws.Range('A73').Value = 65
$ws.Range('B73').Value = '$(F_{p})_{3}$'
$ws.Range('C73').Value = 0.000001115444055315405
$ws.Range('A74').Value = 33
$ws.Range('B74').Value = '$(F_{p})_{0}$'
$ws.Range('C74').Value = 0.000001029312882318999
$ws.Range('A75').Value = 21
$ws.Range('B75').Value = 'h$_{rs}^{3}$'
$ws.Range('C75').Value = 0.000001010561223360137
$ws.Range('A76').Value = 41
$ws.Range('B76').Value = '$(F_{r})_{0}$'
$ws.Range('C76').Value = 0.0000008977462137763174
$ws.Range('A77').Value = 87
$ws.Range('B77').Value = '$(\langle rs\vert rs \rangle)_{1}$'
$ws.Range('C77').Value = 0.0000008345308863706631
$ws.Range('A78').Value = 7
$ws.Range('B78').Value = 'h$_{pq}^{3}$'
$ws.Range('C78').Value = 0.0000006829527940751493
$ws.Range('A79').Value = 4
$ws.Range('B79').Value = 'h$_{pq}^{0}$'
$ws.Range('C79').Value = 0.0000005175440030462878
$ws.Range('A80').Value = 85
$ws.Range('B80').Value = '$(\langle pq \vert pq \rangle)_{1}$'
$ws.Range('C80').Value = 0.0000005012189302508587
$ws.Range('A81').Value = 103
$ws.Range('B81').Value = '$(\langle rs\vert rs \rangle)_{3}$'
$ws.Range('C81').Value = 0.0000004930600012563717
$ws.Range('A82').Value = 83
$ws.Range('B82').Value = '$(\langle pp \vert pp \rangle)_{1}$'
$ws.Range('C82').Value = 0.0000004927169212199049
$ws.Range('A83').Value = 6
$ws.Range('B83').Value = 'h$_{pq}^{2}$'
$ws.Range('C83').Value = 0.0000004732630602851758
$ws.Range('A84').Value = 28
$ws.Range('B84').Value = 'F$_{qs}$'
$ws.Range('C84').Value = 0.0000003603659128482466
$ws.Range('A85').Value = 89
$ws.Range('B85').Value = '$(\langle pq \vert rs \rangle)_{2}$'
$ws.Range('C85').Value = 0.0000003188344659274526
$ws.Range('A86').Value = 20
$ws.Range('B86').Value = 'h$_{rs}^{2}$'
$ws.Range('C86').Value = 0.0000001378075761024632
$ws.Range('A87').Value = 98
$ws.Range('B87').Value = '$(\langle pq \vert sr \rangle)_{3}$'
$ws.Range('C87').Value = 0.000000007042809793413909
$ws.Range('A88').Value = 72
$ws.Range('B88').Value = '$(\langle pq \vert sr \rangle)_{0}$'
$ws.Range('C88').Value = 0.000000005859150272481092
$ws.Range('A89').Value = 82
$ws.Range('B89').Value = '$(\langle pq \vert sr \rangle)_{1}$'
$ws.Range('C89').Value = 0.000000004363823016639468
$ws.Range('A90').Value = 81
$ws.Range('B90').Value = '$(\langle pq \vert rs \rangle)_{1}$'
$ws.Range('C90').Value = 0.000000004282344878950518
$ws.Range('A91').Value = 27
$ws.Range('B91').Value = '$\mathbf{b}$'
$ws.Range('C91').Value = 0.000000003332878082008583
$ws.Range('A92').Value = 40
$ws.Range('B92').Value = '$(\omega_{r})_{0}$'
$ws.Range('C92').Value = 0.000000003218433671846697
$ws.Range('A93').Value = 52
$ws.Range('B93').Value = '$(\omega_{r})_{1}$'
$ws.Range('C93').Value = 0.000000003061761032898048
$ws.Range('A94').Value = 90
$ws.Range('B94').Value = '$(\langle pq \vert sr \rangle)_{2}$'
$ws.Range('C94').Value = 0.000000002749630212097297
$ws.Range('A95').Value = 68
$ws.Range('B95').Value = '$(\omega_{r})_{3}$'
$ws.Range('C95').Value = 0.000000002562520862638928
$ws.Range('A96').Value = 44
$ws.Range('B96').Value = '$\omega_{s}$'
$ws.Range('C96').Value = 0.000000002053595893905851
$ws.Range('A97').Value = 60
$ws.Range('B97').Value = '$(\omega_{r})_{2}$'
$ws.Range('C97').Value = 0.000000001362750414718833
$ws.Range('A98').Value = 36
$ws.Range('B98').Value = '$\omega_{q}$'
$ws.Range('C98').Value = 0.0000000006135002374218353
